$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Outline levels: promote the title and all the "heading-ish"
#    numbered paragraphs to OutlineLevel 1 (-> <w:outlineLvl w:val="0"/>)
# ---------------------------------------------------------------------
$outlineParas = @(1, 3, 6, 15, 18, 21)
foreach ($i in $outlineParas) {
    $d.Paragraphs($i).OutlineLevel = 1
}

# ---------------------------------------------------------------------
# 2) Title paragraph: merge the two runs
#    "Escopo do Projeto (EAP" + ") - Byte INC." into a single run.
# ---------------------------------------------------------------------
$dash = [char]8211
$titleText = "Escopo do Projeto (EAP) " + $dash + " Byte INC."
$titleRange = $d.Paragraphs(1).Range
$titleRange.Find.Execute($titleText, $false, $false, $false, $false, $false, $true, 1, $false, $titleText, 2) | Out-Null

# ---------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark from between the two runs of the
#    "Aprovação e Divulgação" paragraph down into the (previously
#    empty) paragraph that follows it, then merge the
#    "<<O plano..." / ">>" runs into a single run.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$lastRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $lastRange) | Out-Null

$approvalPara = $d.Paragraphs($d.Paragraphs.Count - 1)
$approvalRange = $approvalPara.Range
$eacute = [char]233
$approvalText = "<<O plano " + $eacute + " aprovado aqui depois de pronto.>>"
$approvalRange.Find.Execute($approvalText, $false, $false, $false, $false, $false, $true, 1, $false, $approvalText, 2) | Out-Null

# ---------------------------------------------------------------------
# 4) New styles: "Document Map" paragraph style (MapadoDocumento) and
#    its linked character style (MapadoDocumentoChar).
# ---------------------------------------------------------------------
$mapStyle = $d.Styles.Add("MapadoDocumento", 1)
$mapStyle.NameLocal = "Document Map"
$mapStyle.BaseStyle = $d.Styles("Normal")
$mapStyle.Priority = 99
$mapStyle.UnhideWhenUsed = $true
$mapStyle.Font.Name = "Helvetica"

$mapCharStyle = $d.Styles.Add("MapadoDocumentoChar", 2)
$mapCharStyle.NameLocal = "Mapa do Documento Char"
$mapCharStyle.BaseStyle = $d.Styles("Fontepargpadro")
$mapCharStyle.Priority = 99
$mapCharStyle.Font.Name = "Helvetica"

$mapStyle.LinkStyle = $mapCharStyle
$mapCharStyle.LinkStyle = $mapStyle
